# Update "想去人数" (want-to-go count) values in column F across all four
# sheets, matching the refreshed data snapshot (commit: "Update gh-pages
# to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 5877
$ws.Range("F5").Value  = 5877
$ws.Range("F7").Value  = 2920
$ws.Range("F8").Value  = 1260
$ws.Range("F9").Value  = 391
$ws.Range("F11").Value = 107
$ws.Range("F13").Value = 680
$ws.Range("F14").Value = 170
$ws.Range("F15").Value = 4207
$ws.Range("F16").Value = 4207
$ws.Range("F19").Value = 93
$ws.Range("F21").Value = 190
$ws.Range("F22").Value = 55
$ws.Range("F23").Value = 6313
$ws.Range("F24").Value = 6313
$ws.Range("F28").Value = 429
$ws.Range("F29").Value = 209
$ws.Range("F30").Value = 444
$ws.Range("F31").Value = 5756
$ws.Range("F32").Value = 1610
$ws.Range("F34").Value = 1843
$ws.Range("F35").Value = 5858
$ws.Range("F36").Value = 97
$ws.Range("F38").Value = 82
$ws.Range("F40").Value = 188
$ws.Range("F41").Value = 3943
$ws.Range("F42").Value = 174
$ws.Range("F45").Value = 2387
$ws.Range("F49").Value = 13
$ws.Range("F50").Value = 274
$ws.Range("F52").Value = 15

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 24
$ws.Range("F5").Value = 88

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1403

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1403
$ws.Range("F4").Value  = 5877
$ws.Range("F5").Value  = 5877
$ws.Range("F7").Value  = 2920
$ws.Range("F8").Value  = 1260
$ws.Range("F10").Value = 107
$ws.Range("F13").Value = 680
$ws.Range("F14").Value = 170
$ws.Range("F15").Value = 4207
$ws.Range("F16").Value = 4207
$ws.Range("F19").Value = 93
$ws.Range("F21").Value = 190
$ws.Range("F22").Value = 55
$ws.Range("F23").Value = 6313
$ws.Range("F24").Value = 6313
$ws.Range("F27").Value = 429
$ws.Range("F28").Value = 209
$ws.Range("F29").Value = 88
$ws.Range("F30").Value = 5756
$ws.Range("F31").Value = 1610
$ws.Range("F34").Value = 1843
$ws.Range("F36").Value = 5858
$ws.Range("F37").Value = 97
$ws.Range("F40").Value = 3943
$ws.Range("F45").Value = 2387
$ws.Range("F49").Value = 13
$ws.Range("F50").Value = 274
